$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 holding the 2022-Q3 totals, pushing the
#    existing 2022-Q2 totals row down to row 3.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()

# Bring over A2's number-style (the inserted row has no style of its own
# because column A was empty above it) from the row that got pushed down.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
# B2:D2 picked up a copied style from row 1 on insert; the source data has
# no explicit style on these cells, so clear it back out.
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.54

$wsTotal.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# 2) The existing "2022-Q2" detail sheet becomes the new "2022-Q3" detail
#    sheet (same position / identity); a duplicate of its original content is
#    placed right after it and renamed back to "2022-Q2" so the old figures
#    stay available on their own tab.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)

$wsQ2.Copy($null, $wsQ2)
$wsQ2Copy = $wb.Worksheets.Item(3)

# Rename the original before the copy so the two never collide on name.
$wsQ2.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

# Match page margins used elsewhere in the workbook (Excel's "Normal" preset).
$wsQ2.PageSetup.LeftMargin = 54
$wsQ2.PageSetup.RightMargin = 54
$wsQ2.PageSetup.TopMargin = 72
$wsQ2.PageSetup.BottomMargin = 72
$wsQ2.PageSetup.HeaderMargin = 36
$wsQ2.PageSetup.FooterMargin = 36

# Re-style the header row and the index column to match the bold/centered
# style used on the "总计" sheet.
$wsTotal.Range("B1").Copy()
$wsQ2.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ2.Range("A2:A5").PasteSpecial(-4122)

# Row 1 - headers (unchanged text, kept for clarity).
$wsQ2.Range("B1").Value = "基金代码"
$wsQ2.Range("C1").Value = "基金名称"
$wsQ2.Range("D1").Value = "基金规模"
$wsQ2.Range("E1").Value = "股票总仓位"
$wsQ2.Range("F1").Value = "仓位占比"
$wsQ2.Range("G1").Value = "持有市值(亿元)"
$wsQ2.Range("H1").Value = "仓位排名"

# Fund codes (B) and the D:G figures are stored as text in the source data
# (leading zeros in fund codes must survive).
$wsQ2.Range("B2:B5").NumberFormat = "@"
$wsQ2.Range("D2:G5").NumberFormat = "@"

# Row 2
$wsQ2.Range("A2").Value = 0
$wsQ2.Range("B2").Value = "006679"
$wsQ2.Range("C2").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 A"
$wsQ2.Range("D2").Value = "11.73"
$wsQ2.Range("E2").Value = "93.96"
$wsQ2.Range("F2").Value = "3.07"
$wsQ2.Range("G2").Value = "0.3601"
$wsQ2.Range("H2").Value = 10

# Row 3
$wsQ2.Range("A3").Value = 1
$wsQ2.Range("B3").Value = "162719"
$wsQ2.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
$wsQ2.Range("D3").Value = "11.73"
$wsQ2.Range("E3").Value = "93.96"
$wsQ2.Range("F3").Value = "3.07"
$wsQ2.Range("G3").Value = "0.3601"
$wsQ2.Range("H3").Value = 10

# Row 4
$wsQ2.Range("A4").Value = 2
$wsQ2.Range("B4").Value = "006680"
$wsQ2.Range("C4").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）美元现汇 C"
$wsQ2.Range("D4").Value = "5.92"
$wsQ2.Range("E4").Value = "93.96"
$wsQ2.Range("F4").Value = "3.07"
$wsQ2.Range("G4").Value = "0.1817"
$wsQ2.Range("H4").Value = 10

# Row 5
$wsQ2.Range("A5").Value = 3
$wsQ2.Range("B5").Value = "004243"
$wsQ2.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
$wsQ2.Range("D5").Value = "-11.74"
$wsQ2.Range("E5").Value = "93.96"
$wsQ2.Range("F5").Value = "3.07"
$wsQ2.Range("G5").Value = "-0.3604"
$wsQ2.Range("H5").Value = 10

# The text format was only needed to make Excel store these as text instead
# of re-parsing them into numbers; the source cells carry no explicit style,
# so drop the "@" formatting now that the values are locked in as text.
$wsQ2.Range("B2:B5").ClearFormats()
$wsQ2.Range("D2:G5").ClearFormats()
